$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2333333333333333
$ws.Range("C2").Value = 0.5
$ws.Range("J2").Value = 0.06666666666666667
$ws.Range("P2").Value = 0.1666666666666667
$ws.Range("S2").Value = 0.03333333333333333
$ws.Range("P3").Value = 0.7333333333333333
$ws.Range("S3").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("P5").Value = 1
$ws.Range("O6").Value = 0.1
$ws.Range("Q6").Value = 0.2
$ws.Range("R6").Value = 0.3
$ws.Range("S6").Value = 0.4
$ws.Range("F7").Value = 0.1111111111111111
$ws.Range("J7").Value = 0.1111111111111111
$ws.Range("Q7").Value = 0.3333333333333333
$ws.Range("R7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.3333333333333333
$ws.Range("B8").Value = 0.1153846153846154
$ws.Range("J8").Value = 0.1923076923076923
$ws.Range("Q8").Value = 0.2692307692307692
$ws.Range("R8").Value = 0.1538461538461539
$ws.Range("S8").Value = 0.2692307692307692
$ws.Range("J9").Value = 0.2222222222222222
$ws.Range("Q9").Value = 0.2222222222222222
$ws.Range("R9").Value = 0.3333333333333333
$ws.Range("S9").Value = 0.2222222222222222
$ws.Range("B10").Value = 0.152
$ws.Range("D10").Value = 0.024
$ws.Range("E10").Value = 0.008
$ws.Range("F10").Value = 0.056
$ws.Range("J10").Value = 0.192
$ws.Range("O10").Value = 0.008
$ws.Range("Q10").Value = 0.112
$ws.Range("R10").Value = 0.12
$ws.Range("S10").Value = 0.328
$ws.Range("G11").Value = 0.2
$ws.Range("J11").Value = 0.1333333333333333
$ws.Range("K11").Value = 0.2666666666666667
$ws.Range("L11").Value = 0.4
$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.5
$ws.Range("G13").Value = 0.8
$ws.Range("J13").Value = 0.2
$ws.Range("H15").Value = 0.125
$ws.Range("J15").Value = 0.4375
$ws.Range("K15").Value = 0.0625
$ws.Range("O15").Value = 0.125
$ws.Range("S15").Value = 0.25
$ws.Range("H16").Value = 0.05555555555555555
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("S16").Value = 0.1666666666666667
$ws.Range("H17").Value = 0.1111111111111111
$ws.Range("I17").Value = 0.03703703703703703
$ws.Range("J17").Value = 0.6666666666666666
$ws.Range("K17").Value = 0.03703703703703703
$ws.Range("O17").Value = 0.03703703703703703
$ws.Range("S17").Value = 0.1111111111111111
$ws.Range("H18").Value = 0.1923076923076923
$ws.Range("I18").Value = 0.03846153846153846
$ws.Range("J18").Value = 0.5384615384615384
$ws.Range("K18").Value = 0.03846153846153846
$ws.Range("M18").Value = 0.03846153846153846
$ws.Range("O18").Value = 0.03846153846153846
$ws.Range("S18").Value = 0.1153846153846154
$ws.Range("H19").Value = 0.1666666666666667
$ws.Range("I19").Value = 0.07142857142857142
$ws.Range("J19").Value = 0.3809523809523809
$ws.Range("K19").Value = 0.09523809523809523
$ws.Range("M19").Value = 0.04761904761904762
$ws.Range("O19").Value = 0.09523809523809523
$ws.Range("S19").Value = 0.1428571428571428
